$wb = $excel.ActiveWorkbook
$wsC = $wb.Worksheets.Item("Clientes")
$wsR = $wb.Worksheets.Item("Reservas")

# ---------------------------------------------------------------------------
# 1. Add the 10 new "Reservas" rows (62..71).
# ---------------------------------------------------------------------------

# Column A / B: dates. Row 62 is a literal start date, the remaining rows
# chain off the previous row (+1), except row 67 which restarts from row 62
# (+61) like the existing blocks further up the sheet.
$wsR.Range("A62").Value = 43952
$wsR.Range("B62").Value = 43953

$wsR.Range("A63").Formula = "=A62+1"
$wsR.Range("B63").Formula = "=B62+1"

$wsR.Range("A64").Formula = "=A63+1"
$wsR.Range("B64").Formula = "=B63+1"

$wsR.Range("A65").Formula = "=A64+1"
$wsR.Range("B65").Formula = "=B64+1"

$wsR.Range("A66").Formula = "=A65+1"
$wsR.Range("B66").Formula = "=B65+1"

$wsR.Range("A67").Formula = "=A62+61"
$wsR.Range("B67").Formula = "=B62+61"

$wsR.Range("A68").Formula = "=A67+1"
$wsR.Range("B68").Formula = "=B67+1"

$wsR.Range("A69").Formula = "=A68+1"
$wsR.Range("B69").Formula = "=B68+1"

$wsR.Range("A70").Formula = "=A69+1"
$wsR.Range("B70").Formula = "=B69+1"

$wsR.Range("A71").Formula = "=A70+1"
$wsR.Range("B71").Formula = "=B70+1"

# Apply the same number format the other rows' date cells use (built-in
# "Short Date", numFmtId 14) to the new A/B cells, then propagate that exact
# style to the rest of the block via copy/paste-special so every cell shares
# one cellXf instead of minting a new one each time.
$wsR.Range("A62").NumberFormat = "mm-dd-yy"
$wsR.Range("A62").Copy() | Out-Null
$wsR.Range("A62:B71").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Column C: room type (shared strings already present in the workbook).
$wsR.Range("C62").Value = "single"
$wsR.Range("C63").Value = "twin"
$wsR.Range("C64").Value = "suite"
$wsR.Range("C65").Value = "suite"
$wsR.Range("C66").Value = "superior"
$wsR.Range("C67").Value = "superior"
$wsR.Range("C68").Value = "suite"
$wsR.Range("C69").Value = "suite"
$wsR.Range("C70").Value = "twin"
$wsR.Range("C71").Value = "single"

# Column D: dataReserva = A - 30, same formula pattern/style as existing rows.
$wsR.Range("D62").Formula = "=A62-30"
$wsR.Range("D63").Formula = "=A63-30"
$wsR.Range("D64").Formula = "=A64-30"
$wsR.Range("D65").Formula = "=A65-30"
$wsR.Range("D66").Formula = "=A66-30"
$wsR.Range("D67").Formula = "=A67-30"
$wsR.Range("D68").Formula = "=A68-30"
$wsR.Range("D69").Formula = "=A69-30"
$wsR.Range("D70").Formula = "=A70-30"
$wsR.Range("D71").Formula = "=A71-30"

$wsR.Range("D2").Copy() | Out-Null
$wsR.Range("D62:D71").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Column E: numeroPessoas.
$wsR.Range("E62").Value = 1
$wsR.Range("E63").Value = 2
$wsR.Range("E64").Value = 3
$wsR.Range("E65").Value = 3
$wsR.Range("E66").Value = 4
$wsR.Range("E67").Value = 4
$wsR.Range("E68").Value = 3
$wsR.Range("E69").Value = 3
$wsR.Range("E70").Value = 2
$wsR.Range("E71").Value = 1

# Column F: clienteNIF, pulled from the matching new Clientes rows (6..15).
$wsR.Range("F62").Formula = "=Clientes!B6"
$wsR.Range("F63").Formula = "=Clientes!B7"
$wsR.Range("F64").Formula = "=Clientes!B8"
$wsR.Range("F65").Formula = "=Clientes!B9"
$wsR.Range("F66").Formula = "=Clientes!B10"
$wsR.Range("F67").Formula = "=Clientes!B11"
$wsR.Range("F68").Formula = "=Clientes!B12"
$wsR.Range("F69").Formula = "=Clientes!B13"
$wsR.Range("F70").Formula = "=Clientes!B14"
$wsR.Range("F71").Formula = "=Clientes!B15"

# Column G: the generated INSERT statement.
$insertFormula = '=_xlfn.CONCAT("INSERT INTO reserva (dataEntrada, dataSaida, tipoQuarto, dataReserva, numeroPessoas, clienteNIF) VALUES(''",TEXT({0},"AAAA-MM-DD"),"'', ''",TEXT({1},"AAAA-MM-DD"),"'', ''",{2},"'', ''",TEXT({3},"AAAA-MM-DD"),"'', ",{4},", ",{5},");")'
for ($r = 62; $r -le 71; $r++) {
    $f = $insertFormula -f "A$r", "B$r", "C$r", "D$r", "E$r", "F$r"
    $wsR.Range("G$r").Formula = $f
}

# ---------------------------------------------------------------------------
# 2. Update sheet dimensions / selections.
# ---------------------------------------------------------------------------
$wsC.Range("C40").Select() | Out-Null
$wsR.Activate() | Out-Null
$wsR.Range("C63").Select() | Out-Null
